{"js": "// Replace the date line and every \"A\u00d7B=\" multiplication prompt in the\n// document with its updated value, per the commit's regenerated content.\nconst replacements = [\n  [\"2025-04-13 Sunday\", \"2025-04-14 Monday\"],\n  [\"972\u00d76=\", \"397\u00d74=\"],\n  [\"902\u00d78=\", \"275\u00d72=\"],\n  [\"548\u00d79=\", \"361\u00d74=\"],\n  [\"648\u00d74=\", \"922\u00d79=\"],\n  [\"468\u00d74=\", \"316\u00d77=\"],\n  [\"740\u00d76=\", \"735\u00d76=\"],\n  [\"133\u00d73=\", \"583\u00d78=\"],\n  [\"844\u00d79=\", \"516\u00d72=\"],\n  [\"495\u00d74=\", \"915\u00d73=\"],\n  [\"180\u00d78=\", \"271\u00d75=\"],\n  [\"716\u00d72=\", \"403\u00d72=\"],\n  [\"877\u00d77=\", \"905\u00d73=\"],\n  [\"231\u00d74=\", \"197\u00d73=\"],\n  [\"975\u00d79=\", \"735\u00d74=\"],\n  [\"296\u00d77=\", \"516\u00d76=\"],\n  [\"897\u00d77=\", \"780\u00d79=\"],\n  [\"452\u00d78=\", \"182\u00d72=\"],\n  [\"785\u00d77=\", \"356\u00d72=\"],\n  [\"345\u00d78=\", \"424\u00d74=\"],\n  [\"564\u00d77=\", \"581\u00d75=\"],\n  [\"236\u00d73=\", \"840\u00d75=\"],\n  [\"892\u00d75=\", \"235\u00d76=\"],\n  [\"171\u00d76=\", \"593\u00d75=\"],\n  [\"772\u00d77=\", \"138\u00d78=\"],\n  [\"370\u00d77=\", \"615\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00d7B=\" multiplication prompt in the\n# document with its updated value, per the commit's regenerated content.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-13 Sunday\", \"2025-04-14 Monday\"),\n    @(\"972\u00d76=\", \"397\u00d74=\"),\n    @(\"902\u00d78=\", \"275\u00d72=\"),\n    @(\"548\u00d79=\", \"361\u00d74=\"),\n    @(\"648\u00d74=\", \"922\u00d79=\"),\n    @(\"468\u00d74=\", \"316\u00d77=\"),\n    @(\"740\u00d76=\", \"735\u00d76=\"),\n    @(\"133\u00d73=\", \"583\u00d78=\"),\n    @(\"844\u00d79=\", \"516\u00d72=\"),\n    @(\"495\u00d74=\", \"915\u00d73=\"),\n    @(\"180\u00d78=\", \"271\u00d75=\"),\n    @(\"716\u00d72=\", \"403\u00d72=\"),\n    @(\"877\u00d77=\", \"905\u00d73=\"),\n    @(\"231\u00d74=\", \"197\u00d73=\"),\n    @(\"975\u00d79=\", \"735\u00d74=\"),\n    @(\"296\u00d77=\", \"516\u00d76=\"),\n    @(\"897\u00d77=\", \"780\u00d79=\"),\n    @(\"452\u00d78=\", \"182\u00d72=\"),\n    @(\"785\u00d77=\", \"356\u00d72=\"),\n    @(\"345\u00d78=\", \"424\u00d74=\"),\n    @(\"564\u00d77=\", \"581\u00d75=\"),\n    @(\"236\u00d73=\", \"840\u00d75=\"),\n    @(\"892\u00d75=\", \"235\u00d76=\"),\n    @(\"171\u00d76=\", \"593\u00d75=\"),\n    @(\"772\u00d77=\", \"138\u00d78=\"),\n    @(\"370\u00d77=\", \"615\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
